# "contingencies with rene fine"
# Insert two new line rows (line7, line8) into the lines table, right after
# line6 / before extr1. This pushes extr1..extr8 down by two rows. Also
# flip the in_service flag for extr5 and extr6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at row 8 (pushes old rows 8-15 down to 10-17)
$ws.Rows.Item(8).Resize(2).Insert()

# Give column A of the two new rows the same formatting (bold/border/center)
# used by the rest of the index column.
$ws.Range("A7").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)

# New row 8: line7
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# New row 9: line8
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Renumber column A (plain sequential index, not a formula) for the rows
# that were pushed down by the insertion (old rows 8-15 -> new rows 10-17).
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15

# extr5 (now row 14) and extr6 (now row 15) flip their in_service flag.
$ws.Range("E14").Value = $false
$ws.Range("E15").Value = $true
